# Generate Report for Handoff
# Adds a new localization-status row (cb402c2a-...) alongside the refreshed
# existing row (a9424579-... -> 74cdc530-...) on the Overview, zh-cn and
# de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "a9424579-d46c-422d-ac2e-db9dc5fcbe68"
$guid1   = "74cdc530-baab-48d0-bcfc-7225f2ec2db6"
$guid2   = "cb402c2a-e2a9-434e-b7b0-e1aef770442d"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Refresh existing row (row 2): UUID + generate-date change.
$wsOverview.Range("A2").Value = "$guid1.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$guid1.md"
$wsOverview.Range("G2").Value = "2017-02-21 10:53:49"

# Append the new row (row 3) to the Overview table.
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()

$wsOverview.Range("A3").Value = "$guid2.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2017-02-21 10:53:49"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/7a6f3e3f62068df00a148fdedb9dbc778fbfdea8/e2e/$guid2.md", "", "", "e2e\$guid2.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Refresh existing row (row 2): UUID, handoff file name and timestamp.
$wsZh.Range("A2").Value = "$guid1.md"
$wsZh.Hyperlinks.Item(1).TextToDisplay = "$guid1.md"
$wsZh.Range("G2").Value = "$guid1.9077da9b8686fe5c9b15cdc4b55950110e56b658.zh-cn.xlf"
$wsZh.Range("H2").Value = "2017-02-21 10:53:33"

# Append the new row (row 3) to the zh-cn table.
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = "$guid2.719309e1036a43ec2a0cc7d6164ddac16cfab799.zh-cn.xlf"
$wsZh.Range("H3").Value = "2017-02-21 10:53:33"
$wsZh.Range("L3").Value = "0001-01-01 00:00:00"
$wsZh.Range("O3").Value = "True"
$wsZh.Range("Q3").Value = "True"
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/7a6f3e3f62068df00a148fdedb9dbc778fbfdea8/e2e/$guid2.md", "", "", "$guid2.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Refresh existing row (row 2): UUID, handoff file name and dependency flag.
$wsDe.Range("A2").Value = "$guid1.md"
$wsDe.Hyperlinks.Item(1).TextToDisplay = "$guid1.md"
$wsDe.Range("G2").Value = "$guid1.9077da9b8686fe5c9b15cdc4b55950110e56b658.de-de.xlf"
$wsDe.Range("Q2").Value = "False"

# Append the new row (row 3) to the de-de table.
$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = "$guid2.719309e1036a43ec2a0cc7d6164ddac16cfab799.de-de.xlf"
$wsDe.Range("H3").Value = "2017-02-21 10:53:49"
$wsDe.Range("L3").Value = "0001-01-01 00:00:00"
$wsDe.Range("O3").Value = "True"
$wsDe.Range("Q3").Value = "True"
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/7a6f3e3f62068df00a148fdedb9dbc778fbfdea8/e2e/$guid2.md", "", "", "$guid2.md")
